$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I2: a plain-text copy of the e-mail address already shown in E2
# (same shared string, no hyperlink formatting).
$ws.Range("I2").Value = "anaisabelculajay@gmail.com"

# Turn E2 into a live mailto: hyperlink (adds the relationship) and apply the
# built-in "Hyperlink" cell style, matching the existing G2 hyperlink cell.
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:anaisabelculajay@gmail.com") | Out-Null
$ws.Range("E2").Style = "Hyperlink"

# Leave the selection where the author left it when they saved.
$ws.Range("E4").Select()
